$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.5036099671744125
$ws.Range("C2").Value = 0.04052034902352375
$ws.Range("D2").Value = 0.1878942985749461
$ws.Range("E2").Value = 0.402314294157577
$ws.Range("F2").Value = 3.487510828285707
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.4667827300587248

$ws.Range("B3").Value = 0.4724063787156183
$ws.Range("C3").Value = 0.03569615230428269
$ws.Range("D3").Value = 0.1760989280799379
$ws.Range("E3").Value = 0.3508673324255938
$ws.Range("F3").Value = 3.26612213088967
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.4333145590494496

$ws.Range("B4").Value = 0.4537886395938244
$ws.Range("C4").Value = 0.03278668185521383
$ws.Range("D4").Value = 0.1688137653518993
$ws.Range("E4").Value = 0.3194104324274178
$ws.Range("F4").Value = 3.130498158311099
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.4133131964092911

$ws.Range("B5").Value = 0.4463366259800807
$ws.Range("C5").Value = 0.03161388729603232
$ws.Range("D5").Value = 0.1658337207010021
$ws.Range("E5").Value = 0.3066218381585486
$ws.Range("F5").Value = 3.075301371088642
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.4052987025104358

$ws.Range("B6").Value = 0.4451073337293963
$ws.Range("C6").Value = 0.03141990975724696
$ws.Range("D6").Value = 0.1653381897343706
$ws.Range("E6").Value = 0.3045000542299618
$ws.Range("F6").Value = 3.066140083660372
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.4039760823100949

$ws.Range("B7").Value = 0.453687594733168
$ws.Range("C7").Value = 0.03277081366729817
$ws.Range("D7").Value = 0.1687736217985218
$ws.Range("E7").Value = 0.3192378417100912
$ws.Range("F7").Value = 3.129753476227449
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.4132045607217947

$ws.Range("B8").Value = 0.4927378209911524
$ws.Range("C8").Value = 0.03884580997117837
$ws.Range("D8").Value = 0.1838358066364947
$ws.Range("E8").Value = 0.3845463406985772
$ws.Range("F8").Value = 3.411107596284069
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.455128116654862

$ws.Range("B9").Value = 0.5736739766610981
$ws.Range("C9").Value = 0.05119463413447534
$ws.Range("D9").Value = 0.2130585166596433
$ws.Range("E9").Value = 0.513793844762958
$ws.Range("F9").Value = 3.965631833296072
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 0.5417697655766744

$ws.Range("B10").Value = 0.6358907497860002
$ws.Range("C10").Value = 0.06055955985114281
$ws.Range("D10").Value = 0.2343733135175228
$ws.Range("E10").Value = 0.6096635881757919
$ws.Range("F10").Value = 4.37524365221384
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 0.6082450606347152

$ws.Range("B11").Value = 0.6648141964093668
$ws.Range("C11").Value = 0.0648892788139932
$ws.Range("D11").Value = 0.2440442864441366
$ws.Range("E11").Value = 0.6535182957295973
$ws.Range("F11").Value = 4.562177163855381
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 0.6391254755382363

$ws.Range("B12").Value = 0.6758576760110486
$ws.Range("C12").Value = 0.06653930966228927
$ws.Range("D12").Value = 0.2477034236526094
$ws.Range("E12").Value = 0.6701632973675515
$ws.Range("F12").Value = 4.633058606358816
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 0.6509132733907279

$ws.Range("B13").Value = 0.6734752019106054
$ws.Range("C13").Value = 0.06618347340581465
$ws.Range("D13").Value = 0.2469154893754819
$ws.Range("E13").Value = 0.6665767486053653
$ws.Range("F13").Value = 4.617788725657647
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 0.6483703454627516

$ws.Range("B14").Value = 0.6657209214454838
$ws.Range("C14").Value = 0.06502481555834549
$ws.Range("D14").Value = 0.2443453836619369
$ws.Range("E14").Value = 0.6548869062502405
$ws.Range("F14").Value = 4.568006699117689
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 0.6400933683956112

$ws.Range("B15").Value = 0.6609830716621161
$ws.Range("C15").Value = 0.06431648010985214
$ws.Range("D15").Value = 0.2427707397322649
$ws.Range("E15").Value = 0.6477316085807558
$ws.Range("F15").Value = 4.537526212058935
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 0.6350357885032736

$ws.Range("B16").Value = 0.6340131557639666
$ws.Range("C16").Value = 0.06027804122557257
$ws.Range("D16").Value = 0.233740834203445
$ws.Range("E16").Value = 0.6068027451735958
$ws.Range("F16").Value = 4.363039811763372
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 0.6062400069798457

$ws.Range("B17").Value = 0.61762805129149
$ws.Range("C17").Value = 0.05781872960021417
$ws.Range("D17").Value = 0.2281952259479567
$ws.Range("E17").Value = 0.5817587434397495
$ws.Range("F17").Value = 4.256156879186477
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 0.5887401450713128

$ws.Range("B18").Value = 0.6082621038130753
$ws.Range("C18").Value = 0.05641071874924819
$ws.Range("D18").Value = 0.2250031285518759
$ws.Range("E18").Value = 0.5673767794070841
$ws.Range("F18").Value = 4.19473656342538
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 0.5787348517136195

$ws.Range("B19").Value = 0.6051009270976806
$ws.Range("C19").Value = 0.0559350965476284
$ws.Range("D19").Value = 0.2239219077478225
$ws.Range("E19").Value = 0.562511092318573
$ws.Range("F19").Value = 4.173950073580272
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 0.5753575057074443

$ws.Range("B20").Value = 0.6193662255465711
$ws.Range("C20").Value = 0.05807984943200495
$ws.Range("D20").Value = 0.2287858106899421
$ws.Range("E20").Value = 0.5844223494057417
$ws.Range("F20").Value = 4.267528907387884
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 0.5905967944376584

$ws.Range("B21").Value = 0.667996064900592
$ws.Range("C21").Value = 0.06536485388213009
$ws.Range("D21").Value = 0.2451003641552916
$ws.Range("E21").Value = 0.6583194342497478
$ws.Range("F21").Value = 4.582626280412626
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 0.6425219487595939

$ws.Range("B22").Value = 0.7003084580388759
$ws.Range("C22").Value = 0.07018721421148655
$ws.Range("D22").Value = 0.2557453611648839
$ws.Range("E22").Value = 0.7068395096332551
$ws.Range("F22").Value = 4.789110879957661
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 0.6770072031177108

$ws.Range("B23").Value = 0.6830137140461261
$ws.Range("C23").Value = 0.06760767771635301
$ws.Range("D23").Value = 0.2500653350336393
$ws.Range("E23").Value = 0.6809218496353111
$ws.Range("F23").Value = 4.678853257815661
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 0.6585508671337834

$ws.Range("B24").Value = 0.6185802286740341
$ws.Range("C24").Value = 0.05796177887495446
$ws.Range("D24").Value = 0.2285188193744148
$ws.Range("E24").Value = 0.5832180830372806
$ws.Range("F24").Value = 4.262387526469809
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 0.5897572301927312

$ws.Range("B25").Value = 0.5513021923533188
$ws.Range("C25").Value = 0.04780457923702386
$ws.Range("D25").Value = 0.205182891130022
$ws.Range("E25").Value = 0.4786830788613656
$ws.Range("F25").Value = 3.815268189410062
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.5178445011083284
